$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 112083128
$ws.Range("B7").Value = 77186
$ws.Range("C7").Value = 'Ovaliderad'
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 353
$ws.Range("F7").Value = 'Dvärgbägarlav'
$ws.Range("G7").Value = 'Cladonia parasitica'
$ws.Range("H7").Value = '(Hoffm.) Hoffm.'
$ws.Range("P7").Value = 'Sörskog Skallberget, Vrm'
$ws.Range("Q7").Value = 413190.1061828797
$ws.Range("R7").Value = 6656475.01450387
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = 'Värmland'
$ws.Range("U7").Value = 'Hagfors'
$ws.Range("V7").Value = 'Värmland'
$ws.Range("W7").Value = 'Ekshärad'
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = '2023-09-11'
$ws.Range("Z7").Value = '00:00'
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = '2023-09-11'
$ws.Range("AB7").Value = '00:00'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = 'anders tedeholm'
$ws.Range("AX7").Value = 'anders tedeholm'
$ws.Range("Y7").ClearFormats()
$ws.Range("AA7").ClearFormats()

# Row 8
$ws.Range("A8").Value = 112083126
$ws.Range("B8").Value = 78536
$ws.Range("C8").Value = 'Ovaliderad'
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 229497
$ws.Range("F8").Value = 'Korallblylav'
$ws.Range("G8").Value = 'Parmeliella triptophylla'
$ws.Range("H8").Value = '(Ach.) Müll.Arg.'
$ws.Range("P8").Value = 'Sörskog Skallberget, Vrm'
$ws.Range("Q8").Value = 413016.7201701452
$ws.Range("R8").Value = 6656341.641577623
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 'Värmland'
$ws.Range("U8").Value = 'Hagfors'
$ws.Range("V8").Value = 'Värmland'
$ws.Range("W8").Value = 'Ekshärad'
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = '2023-09-11'
$ws.Range("Z8").Value = '00:00'
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = '2023-09-11'
$ws.Range("AB8").Value = '00:00'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = 'anders tedeholm'
$ws.Range("AX8").Value = 'anders tedeholm'
$ws.Range("Y8").ClearFormats()
$ws.Range("AA8").ClearFormats()

# Row 9
$ws.Range("A9").Value = 112083112
$ws.Range("B9").Value = 79444
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 1049
$ws.Range("F9").Value = 'Kortskaftad ärgspik'
$ws.Range("G9").Value = 'Microcalicium ahlneri'
$ws.Range("H9").Value = 'Tibell'
$ws.Range("P9").Value = 'Sörskog Skallberget, Vrm'
$ws.Range("Q9").Value = 412283.7604491137
$ws.Range("R9").Value = 6656072.080045181
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Värmland'
$ws.Range("U9").Value = 'Hagfors'
$ws.Range("V9").Value = 'Värmland'
$ws.Range("W9").Value = 'Ekshärad'
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = '2023-09-11'
$ws.Range("Z9").Value = '00:00'
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = '2023-09-11'
$ws.Range("AB9").Value = '00:00'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = 'anders tedeholm'
$ws.Range("AX9").Value = 'anders tedeholm'
$ws.Range("Y9").ClearFormats()
$ws.Range("AA9").ClearFormats()

# Row 10
$ws.Range("A10").Value = 112083118
$ws.Range("B10").Value = 94134
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 53
$ws.Range("F10").Value = 'Vedtrappmossa'
$ws.Range("G10").Value = 'Crossocalyx hellerianus'
$ws.Range("H10").Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range("P10").Value = 'Sörskog Skallberget, Vrm'
$ws.Range("Q10").Value = 412576.6879626553
$ws.Range("R10").Value = 6656303.56951345
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Värmland'
$ws.Range("U10").Value = 'Hagfors'
$ws.Range("V10").Value = 'Värmland'
$ws.Range("W10").Value = 'Ekshärad'
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = '2023-09-11'
$ws.Range("Z10").Value = '00:00'
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = '2023-09-11'
$ws.Range("AB10").Value = '00:00'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = 'anders tedeholm'
$ws.Range("AX10").Value = 'anders tedeholm'
$ws.Range("Y10").ClearFormats()
$ws.Range("AA10").ClearFormats()

# Row 11
$ws.Range("A11").Value = 112083125
$ws.Range("B11").Value = 89369
$ws.Range("C11").Value = 'Ovaliderad'
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 5447
$ws.Range("F11").Value = 'Vedticka'
$ws.Range("G11").Value = 'Fuscoporia viticola'
$ws.Range("H11").Value = '(Schwein.) Murrill'
$ws.Range("P11").Value = 'Sörskog Skallberget, Vrm'
$ws.Range("Q11").Value = 413015.9403039298
$ws.Range("R11").Value = 6656414.640994807
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Värmland'
$ws.Range("U11").Value = 'Hagfors'
$ws.Range("V11").Value = 'Värmland'
$ws.Range("W11").Value = 'Ekshärad'
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = '2023-09-11'
$ws.Range("Z11").Value = '00:00'
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = '2023-09-11'
$ws.Range("AB11").Value = '00:00'
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AW11").Value = 'anders tedeholm'
$ws.Range("AX11").Value = 'anders tedeholm'
$ws.Range("Y11").ClearFormats()
$ws.Range("AA11").ClearFormats()

# Row 12
$ws.Range("A12").Value = 112083110
$ws.Range("B12").Value = 78107
$ws.Range("C12").Value = 'Ovaliderad'
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 6453
$ws.Range("F12").Value = 'Vedskivlav'
$ws.Range("G12").Value = 'Hertelidea botryosa'
$ws.Range("H12").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("P12").Value = 'Sörskog Skallberget, Vrm'
$ws.Range("Q12").Value = 412205.6393663768
$ws.Range("R12").Value = 6656050.944565876
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 'Värmland'
$ws.Range("U12").Value = 'Hagfors'
$ws.Range("V12").Value = 'Värmland'
$ws.Range("W12").Value = 'Ekshärad'
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = '2023-09-11'
$ws.Range("Z12").Value = '00:00'
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = '2023-09-11'
$ws.Range("AB12").Value = '00:00'
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AW12").Value = 'anders tedeholm'
$ws.Range("AX12").Value = 'anders tedeholm'
$ws.Range("Y12").ClearFormats()
$ws.Range("AA12").ClearFormats()

# Row 13
$ws.Range("A13").Value = 112083127
$ws.Range("B13").Value = 77604
$ws.Range("C13").Value = 'Ovaliderad'
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 6450
$ws.Range("F13").Value = 'Skuggblåslav'
$ws.Range("G13").Value = 'Hypogymnia vittata'
$ws.Range("H13").Value = '(Ach.) Parrique'
$ws.Range("P13").Value = 'Sörskog Skallberget, Vrm'
$ws.Range("Q13").Value = 413051.8096683071
$ws.Range("R13").Value = 6656343.312587639
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = 'Värmland'
$ws.Range("U13").Value = 'Hagfors'
$ws.Range("V13").Value = 'Värmland'
$ws.Range("W13").Value = 'Ekshärad'
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = '2023-09-11'
$ws.Range("Z13").Value = '00:00'
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = '2023-09-11'
$ws.Range("AB13").Value = '00:00'
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AW13").Value = 'anders tedeholm'
$ws.Range("AX13").Value = 'anders tedeholm'
$ws.Range("Y13").ClearFormats()
$ws.Range("AA13").ClearFormats()

# Row 14
$ws.Range("A14").Value = 112083111
$ws.Range("B14").Value = 90666
$ws.Range("C14").Value = 'Ovaliderad'
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 4364
$ws.Range("F14").Value = 'Dropptaggsvamp'
$ws.Range("G14").Value = 'Hydnellum ferrugineum'
$ws.Range("H14").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("P14").Value = 'Sörskog Skallberget, Vrm'
$ws.Range("Q14").Value = 412204.6634863199
$ws.Range("R14").Value = 6655988.977203708
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 'Värmland'
$ws.Range("U14").Value = 'Hagfors'
$ws.Range("V14").Value = 'Värmland'
$ws.Range("W14").Value = 'Ekshärad'
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = '2023-09-11'
$ws.Range("Z14").Value = '00:00'
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = '2023-09-11'
$ws.Range("AB14").Value = '00:00'
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AW14").Value = 'anders tedeholm'
$ws.Range("AX14").Value = 'anders tedeholm'
$ws.Range("Y14").ClearFormats()
$ws.Range("AA14").ClearFormats()
